# This workbook stores the Price column (D) as plain text (inline strings),
# even though most values look numeric (e.g. "205.95"). Excel's COM layer
# auto-converts such text into numbers when assigned directly, so for any
# new Price value that would parse as a plain number we first mark the
# cell's number format as Text ("@") to force the value to stay textual,
# matching the original authoring. Values that are unambiguously
# non-numeric (contain two dots, thousands separators, or the special
# subscript-digit glyphs) do not need this and are set directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $value
}

# Row 2 - Bitcoin
Set-TextValue "D2" "26.845.86"
$ws.Range("E2").Value = "  -0.08%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.543.46"
$ws.Range("E3").Value = "  -1.48%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.23%  "

# Row 5 - BNB
Set-TextValue "D5" "205.95"
$ws.Range("E5").Value = "  -0.32%  "

# Row 6 - XRP
Set-TextValue "D6" "0.486"
$ws.Range("E6").Value = "  -1.04%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.23%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  -0.58%  "

# Row 9 - Solana
Set-TextValue "D9" "21.39"
$ws.Range("E9").Value = "  -2.81%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.72%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -1.11%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "1.763.00"
$ws.Range("E12").Value = "  -1.52%  "

# Row 13 - WrappedEther
Set-TextValue "D13" "1.543.52"
$ws.Range("E13").Value = "  -1.37%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -1.42%  "

# Row 15 - Polygon
Set-TextValue "D15" "0.509"
$ws.Range("E15").Value = "  -1.02%  "

# Row 16 - WrappedBTC
Set-TextValue "D16" "26.845.13"
$ws.Range("E16").Value = "  -0.06%  "

# Row 17 - Litecoin
Set-TextValue "D17" "61.29"
$ws.Range("E17").Value = "  -0.31%  "

# Row 18 - BitcoinCash
Set-TextValue "D18" "215.27"
$ws.Range("E18").Value = "  +0.14%  "

# Rows 19 & 20 swap: ShibaInu <-> Chainlink
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D19" "7.22"
$ws.Range("E19").Value = "  -2.42%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = [string]::Concat("0.0", [string][char]0x2083, "0682")
$ws.Range("E20").Value = "  +0.62%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +0.24%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -2.51%  "

# Row 23 - Avalanche
Set-TextValue "D23" "9.16"
$ws.Range("E23").Value = "  -1.71%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -2.55%  "

# Row 25 - Monero
Set-TextValue "D25" "152.87"
$ws.Range("E25").Value = "  -0.77%  "

# Row 26 - Cosmos
Set-TextValue "D26" "6.62"
$ws.Range("E26").Value = "  -1.62%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "14.85"
$ws.Range("E27").Value = "  -0.81%  "

# Row 28 - BinanceUSD
$ws.Range("E28").Value = "  +0.24%  "

# Row 29 - Stellar
$ws.Range("E29").Value = "  -0.67%  "

# Rows 30 & 31 swap: Hedera <-> PancakeSwap
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D30" "1.10"
$ws.Range("E30").Value = "  -1.44%  "

$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D31" "0.0458"
$ws.Range("E31").Value = "  -1.93%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +1.30%  "

# Row 33 - Maker
Set-TextValue "D33" "1.366.35"
$ws.Range("E33").Value = "  -2.63%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  +0.28%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  -1.11%  "

# Row 36 - TrustWalletToken
Set-TextValue "D36" "0.957"
$ws.Range("E36").Value = "  +2.25%  "

# Row 37 - HuobiToken
$ws.Range("E37").Value = "  +0.06%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  +1.16%  "

# Row 39 - ImmutableX
Set-TextValue "D39" "0.521"
$ws.Range("E39").Value = "  -1.16%  "

# Row 40 - ARBITRUM
$ws.Range("E40").Value = "  -1.10%  "

# Rows 41 & 42 swap: FraxShare <-> PaxDollar
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D41" "1.00"
$ws.Range("E41").Value = "  +0.23%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D42" "5.71"
$ws.Range("E42").Value = "  +7.38%  "

# Row 43 - WEMIXToken
Set-TextValue "D43" "0.990"
$ws.Range("E43").Value = "  +0.01%  "

# Row 44 - MXToken
$ws.Range("E44").Value = "  +1.64%  "

# Row 45 - Aave
Set-TextValue "D45" "63.17"
$ws.Range("E45").Value = "  -0.15%  "

# Row 46 - RenderToken
Set-TextValue "D46" "1.72"
$ws.Range("E46").Value = "  -3.78%  "

# Row 47 - RocketPoolETH
Set-TextValue "D47" "1.677.53"
$ws.Range("E47").Value = "  -1.50%  "

# Row 48 - Quant
Set-TextValue "D48" "84.23"
$ws.Range("E48").Value = "  -2.22%  "

# Row 49 - Cronos
Set-TextValue "D49" "0.0514"
$ws.Range("E49").Value = "  +4.39%  "

# Row 50 - BabyDogeCoin
$ws.Range("D50").Value = [string]::Concat("0.0", [string][char]0x2087, "0967")
$ws.Range("E50").Value = "  -1.76%  "

# Row 51 - Algorand -> USDD
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
Set-TextValue "D51" "1.00"
$ws.Range("E51").Value = "  +0.28%  "
